$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 627
$ws.Range("F9").Value = 6362
$ws.Range("F16").Value = 2627
$ws.Range("F35").Value = 248
$ws.Range("F41").Value = 19
$ws.Range("F43").Value = 2284
$ws.Range("F46").Value = 132

# 演出 (Performance) sheet
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F24").Value = 380

# 本地生活 (Local Life) sheet
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1502
$ws.Range("F10").Value = 2528
$ws.Range("F11").Value = 856
$ws.Range("F13").Value = 30

# 全部类型 (All Types) sheet
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 627
$ws.Range("F9").Value = 1502
$ws.Range("F11").Value = 2528
$ws.Range("F12").Value = 6362
$ws.Range("F13").Value = 856
$ws.Range("F17").Value = 2627
$ws.Range("F35").Value = 248
$ws.Range("F39").Value = 19
$ws.Range("F44").Value = 2284
$ws.Range("F46").Value = 132
